$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '90.961.77'
$ws.Range("E2").Value = '  +1.67%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.174.59'
$ws.Range("E3").Value = '  +4.53%  '

$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.84'
$ws.Range("E5").Value = '  +2.26%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '629.41'
$ws.Range("E6").Value = '  +2.79%  '

$ws.Range("E7").Value = '  +30.32%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.374'
$ws.Range("E8").Value = '  +4.16%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  -0.06%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.170.95'
$ws.Range("E10").Value = '  +4.49%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.752'
$ws.Range("E11").Value = '  +12.28%  '

$ws.Range("E12").Value = '  +8.77%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000247'
$ws.Range("E13").Value = '  +3.64%  '

$ws.Range("E14").Value = '  +6.30%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '35.13'
$ws.Range("E15").Value = '  +9.04%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '90.570.32'
$ws.Range("E16").Value = '  +1.46%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.753.85'
$ws.Range("E17").Value = '  +4.54%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.158.23'
$ws.Range("E18").Value = '  +3.81%  '

$ws.Range("E19").Value = '  +11.74%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.47'
$ws.Range("E20").Value = '  +8.32%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '467.78'
$ws.Range("E21").Value = '  +10.50%  '

$ws.Range("E22").Value = '  -3.76%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.13'
$ws.Range("E23").Value = '  +11.15%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.28'
$ws.Range("E24").Value = '  +5.53%  '

$ws.Range("B25").Value = 'NEARProtocol'
$ws.Range("C25").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.88'
$ws.Range("E25").Value = '  +9.59%  '

$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '93.89'
$ws.Range("E26").Value = '  +12.28%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.26'
$ws.Range("E27").Value = '  +5.91%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.326.19'
$ws.Range("E28").Value = '  +3.81%  '

$ws.Range("E29").Value = '  +0.03%  '

$ws.Range("B30").Value = 'Cronos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.162'
$ws.Range("E30").Value = '  +0.80%  '

$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.25'
$ws.Range("E31").Value = '  +10.96%  '

$ws.Range("B32").Value = 'Binance-PegBSC-USD'
$ws.Range("C32").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("E32").Value = '  +0.09%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '28.19'
$ws.Range("E33").Value = '  +24.09%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '525.49'
$ws.Range("E34").Value = '  +4.84%  '

$ws.Range("E35").Value = '  +36.64%  '

$ws.Range("B36").Value = 'PancakeSwap'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.94'
$ws.Range("E36").Value = '  +7.76%  '

$ws.Range("B37").Value = 'dogwifhat'
$ws.Range("C37").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.64'
$ws.Range("E37").Value = '  -2.63%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.96'
$ws.Range("E38").Value = '  +5.09%  '

$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.142'
$ws.Range("E39").Value = '  +7.83%  '

$ws.Range("B40").Value = 'Fetch.AI'
$ws.Range("C40").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.31'
$ws.Range("E40").Value = '  +5.93%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '22.22'
$ws.Range("E41").Value = '  -0.04%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0850'
$ws.Range("E42").Value = '  +24.36%  '

$ws.Range("E43").Value = '  +0.09%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.413'
$ws.Range("E44").Value = '  +14.54%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.98'
$ws.Range("E45").Value = '  +8.25%  '

$ws.Range("E46").Value = '  -0.02%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '150.31'
$ws.Range("E47").Value = '  +1.89%  '

$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.697'
$ws.Range("E48").Value = '  +18.69%  '

$ws.Range("B49").Value = 'OKB'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '45.40'
$ws.Range("E49").Value = '  +4.62%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.36'
$ws.Range("E50").Value = '  +11.79%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.52'
$ws.Range("E51").Value = '  +7.87%  '
